$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary" - refresh the headline metrics for the newly-closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.51   # Current Capital
$summary.Range("B4").Value = -2.49    # Total P&L $
$summary.Range("B5").Value = -1.16    # Total P&L %
$summary.Range("B6").Value = 43        # Total Trades
$summary.Range("B8").Value = 23        # Losing Trades
$summary.Range("B9").Value = 39.53    # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status" - MarketMaking row (row 4) picks up the same trade
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.51000000000001  # Capital
$status.Range("D4").Value = 43                  # Trades
$status.Range("E4").Value = -2.49               # P&L $
$status.Range("F4").Value = -2.49               # P&L %
$status.Range("G4").Value = 39.53               # Win Rate %

# ---------------------------------------------------------------------------
# Helper that appends trade #43 (row 44) to a trade-log style sheet. Both
# "All Trades" and "MarketMaking" get the identical new record.
# ---------------------------------------------------------------------------
function Add-Trade44($ws) {
    $ws.Cells.Item(44, 1).Value = 43

    # Dates/times are stored as literal text in this workbook (column A is
    # numeric trade id, columns B/C are plain strings) - use a leading
    # quote so the engine keeps them as text instead of parsing "2026-02-17"
    # into a date serial number, then drop back to the Normal style so no
    # stray quote-prefix formatting is left behind on the cell.
    $ws.Cells.Item(44, 2).Value = "'2026-02-17"
    $ws.Cells.Item(44, 2).Style = "Normal"
    $ws.Cells.Item(44, 3).Value = "13:27:32"

    $ws.Cells.Item(44, 4).Value = "MarketMaking"
    $ws.Cells.Item(44, 5).Value = "DOWN"
    $ws.Cells.Item(44, 6).Value = 0.29
    $ws.Cells.Item(44, 7).Value = 0.19
    $ws.Cells.Item(44, 8).Value = "CLOSED"
    $ws.Cells.Item(44, 9).Value = -34.4828
    $ws.Cells.Item(44, 10).Value = -0.1
    $ws.Cells.Item(44, 11).Value = 97.51000000000001
    $ws.Cells.Item(44, 12).Value = 0
    $ws.Cells.Item(44, 13).Value = 0
    $ws.Cells.Item(44, 14).Value = 0.6
    $ws.Cells.Item(44, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(44, 16).Value = "early_exit"
    $ws.Cells.Item(44, 17).Value = 0.13
}

Add-Trade44 $wb.Worksheets.Item("All Trades")
Add-Trade44 $wb.Worksheets.Item("MarketMaking")
